# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the existing header row (AC1's style: bold,
# centered, bordered) by copying the format from the last existing header.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Team record is constant for every player row (2-43).
$ws.Range("AD2:AD43").Value = 103
$ws.Range("AE2:AE43").Value = 59
$ws.Range("AF2:AF43").Value = 0
